$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 70
$ws.Range("I5").Value = 70
$ws.Range("K5").Value = 70
$ws.Range("M5").Value = 45
$ws.Range("H10").Value = 5752
$ws.Range("I10").Value = 504
$ws.Range("J10").Value = 11000
$ws.Range("K10").Value = 504
$ws.Range("L10").Value = 11000
$ws.Range("M10").Value = -211
$ws.Range("N10").Value = -11586
$ws.Range("H29").Value = 9900
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 9900
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 29700
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -30262
$ws.Range("H32").Value = 774
$ws.Range("J32").Value = 823.75
$ws.Range("L32").Value = 823.75
$ws.Range("N32").Value = -1475.75
$ws.Range("H33").Value = 367.63635
$ws.Range("I33").Value = 317.875
$ws.Range("K33").Value = 317.875
$ws.Range("M33").Value = -88.875
$ws.Range("H39").Value = 62.57143
$ws.Range("I39").Value = 47.6
$ws.Range("K39").Value = 142.8
$ws.Range("M39").Value = 153.2
$ws.Range("H58").Value = 826
$ws.Range("I58").Value = 826
$ws.Range("K58").Value = 2478
$ws.Range("M58").Value = -2328
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2070.6667
$ws.Range("I122").Value = 2099
$ws.Range("J122").Value = 2014
$ws.Range("K122").Value = 6297
$ws.Range("L122").Value = 6042
$ws.Range("M122").Value = -3847
$ws.Range("N122").Value = -10942

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1869.5555
$ws.Range("I86").Value = 1800.6666
$ws.Range("K86").Value = 1800.6666
$ws.Range("M86").Value = -677.6666
$ws.Range("H89").Value = 1869.5555
$ws.Range("I89").Value = 1800.6666
$ws.Range("K89").Value = 9003.333000000001
$ws.Range("M89").Value = -3387.333000000001

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H82").Value = 20000
$ws.Range("I82").Value = 20000
$ws.Range("K82").Value = 20000
$ws.Range("M82").Value = -19639
$ws.Range("H85").Value = 20000
$ws.Range("I85").Value = 20000
$ws.Range("K85").Value = 20000
$ws.Range("M85").Value = -18752
$ws.Range("H86").Value = 3942.5
$ws.Range("I86").Value = 3940
$ws.Range("K86").Value = 3940
$ws.Range("M86").Value = -2817
$ws.Range("H89").Value = 3942.5
$ws.Range("I89").Value = 3940
$ws.Range("K89").Value = 19700
$ws.Range("M89").Value = -14084

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 73.333336
$ws.Range("J2").Value = 10
$ws.Range("L2").Value = 60
$ws.Range("N2").Value = -286
$ws.Range("H4").Value = 110.1
$ws.Range("I4").Value = 97
$ws.Range("J4").Value = 140.66667
$ws.Range("K4").Value = 291
$ws.Range("L4").Value = 422.00001
$ws.Range("M4").Value = -179
$ws.Range("N4").Value = -646.00001
$ws.Range("H7").Value = 66.333336
$ws.Range("I7").Value = 49.5
$ws.Range("K7").Value = 148.5
$ws.Range("M7").Value = -36.5
$ws.Range("H13").Value = 500
$ws.Range("I13").Value = 700
$ws.Range("J13").Value = 300
$ws.Range("K13").Value = 2100
$ws.Range("L13").Value = 900
$ws.Range("M13").Value = -1932
$ws.Range("N13").Value = -1236
$ws.Range("H23").Value = 396.8
$ws.Range("I23").Value = 149.5
$ws.Range("K23").Value = 448.5
$ws.Range("M23").Value = -213.5
$ws.Range("H34").Value = 533.3333
$ws.Range("H39").Value = 500
$ws.Range("I39").Value = 500
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 1500
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = -1206
$ws.Range("N39").ClearContents()
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("N55").ClearContents()
$ws.Range("H138").Value = 2000
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 2000
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 6000
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -16280
$ws.Range("H140").Value = 1792.8334
$ws.Range("I140").Value = 1792.8334
$ws.Range("K140").Value = 5378.5002
$ws.Range("M140").Value = -198.5002000000004

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 3617.875
$ws.Range("I43").Value = 3617.875
$ws.Range("K43").Value = 3617.875
$ws.Range("M43").Value = -3466.875
$ws.Range("H122").Value = 349.5
$ws.Range("I122").Value = 349.5
$ws.Range("K122").Value = 1048.5
$ws.Range("M122").Value = 1401.5

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1573.75
$ws.Range("I7").Value = 1573.75
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 1573.75
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -1461.75
$ws.Range("N7").ClearContents()
$ws.Range("H22").Value = 377.14285
$ws.Range("I22").Value = 316.66666
$ws.Range("K22").Value = 316.66666
$ws.Range("M22").Value = -21.66665999999998
$ws.Range("H27").Value = 377.14285
$ws.Range("I27").Value = 316.66666
$ws.Range("K27").Value = 316.66666
$ws.Range("M27").Value = -209.66666
$ws.Range("H35").Value = 2396.4285
$ws.Range("I35").Value = 1141
$ws.Range("J35").Value = 5535
$ws.Range("K35").Value = 1141
$ws.Range("L35").Value = 5535
$ws.Range("M35").Value = -805
$ws.Range("N35").Value = -6207
$ws.Range("H61").Value = 5048.8887
$ws.Range("I61").Value = 4234.4287
$ws.Range("K61").Value = 4234.4287
$ws.Range("M61").Value = -4032.4287
$ws.Range("H113").Value = 5048.8887
$ws.Range("I113").Value = 4234.4287
$ws.Range("K113").Value = 4234.4287
$ws.Range("M113").Value = -2064.4287
$ws.Range("H126").Value = 1573.75
$ws.Range("I126").Value = 1573.75
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 4721.25
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -2251.25
$ws.Range("N126").ClearContents()

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 39403.332
$ws.Range("J70").Value = 39403.332
$ws.Range("L70").Value = 39403.332
$ws.Range("N70").Value = -40033.332
$ws.Range("H73").Value = 39403.332
$ws.Range("J73").Value = 39403.332
$ws.Range("L73").Value = 39403.332
$ws.Range("N73").Value = -41587.332
